$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13, shifting rows 13:24 down to 14:25 ---
$ws.Rows("13:13").Insert()

# The insert carries A12's style down into the new A13; that cell must stay
# empty (row 13 in the target only has B13/C13 populated), so clear it.
$ws.Range("A13").Clear()

# --- Row 10 (Objetivos:) — replace the B/C text with the real objectives copy ---
$objetivos = "Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais abordando aspectos bioquímicos importantes na produção de alimentos, e importantes metabólitos. Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna, incluindo as biorefinarias de lignocelulósicos.Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos na resolução de problemas industriais envolvendo processos bioquímicos."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Row 13 (new, under Docentes responsáveis:) — the professor's name moves here ---
$docente = "1814052 - Silvio Silverio da Silva"
$ws.Range("B13").Value = $docente
$ws.Range("C13").Value = $docente

# The freshly-inserted B13 picked up column A's style by default (the <cols>
# range still spans A:B); restore the normal column-B (wrap-text) style by
# copying formatting only from an existing column-B data cell.
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = $docente

# --- Row 14 (Programa resumido:) — replace "Semestral" with the real summary ---
$programaResumido = "Introdução aos processos bioquímicos industriais que incluem o processamento de alimentos, e importantes metabólitos, a manufatura de soros e vacinas, e os conceitos modernos de bioenergia e biorrefinarias."
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido
$ws.Rows("14:14").RowHeight = 60

# --- Row 16 (Programa:) — replace stray date with the actual Portuguese syllabus ---
$programa = "1. Introdução ao Processamento de alimentos: tipos de indústria de alimentos, matériasprimas,fases doprocessamento de produtos alimentícios, conservação/alterações de alimentos, microbiologia dealimentos, alterações bioquímicas em alimentos (oxidação de lipídeos, antioxidantes, escurecimentoenzimático e não enzimático), aflatoxinas, conservantes químicos, toxicantes naturais.2. Discussão e apresentação sobre aspectos bioquímicos importantes na produção de metabólitos por microrganismos e  estudo de casos .3. Manufatura de soros e vacinas Métodosindustriais para a produção de soros e vacinas 4.Biotecnologia de lignocelulósicos: Separação e fermentação das frações e principais processosbioquímicos envolvendo materiais lignocelulósicos.5. Bioenergia, biocombustíveis e biorrefinarias."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Row 19 (Método:) — fill in the evaluation method text ---
$metodo = "A avaliação será feita por meio de prova escrita (P1) e trabalhos (P2)."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# --- Row 20 (Critério:) — final grade formula ---
$criterio = "A nota final (NF) será calculada da seguintes maneira: NF=(P1 + P2)/2"
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# --- Row 21 (Norma de recuperação:) — recovery grade formula ---
$norma = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR)será calculada como MR=(NF+PR)/2"
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

# --- Row 22 (Bibliografia:) — bibliography list, now populated ---
$biblio = "1. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial - Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001.`n2. DEMAIN, A.L., SOLOMON, N.A. (Eds). Manual of industrial microbiology and biotechnology, Washington: American Society for Microbiology, 1986.`n3. WANG, D.C. et al. Fermentation and Enzyme Technology, New York: Wiley-Interscience, 1979.`n4. GAVA, A.J. Princípios de Tecnologia de Alimentos, São Paulo: Nobel, 1983.`n5. LIMA , U. A et al. Biotecnología Industrial, Biotecnologia na produção de alimentos - Série Biotecnología, vol4. Ed. Edgard Blucher,Ltda , 2001.`n6. EVANGELISTA, J. Tecnologia de Alimentos, Rio de Janeiro: Livraria Atheneu, 1987.`n7. CAMARGO R. et al., Tecnologia de produtos Agropecuários- Alimentos, São Paulo: Livraria Nobel, 1984."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio
$ws.Rows("22:22").RowHeight = 120

# --- Row heights for rows that shifted but need their customHeight adjusted ---
$ws.Rows("15:15").RowHeight = 60
$ws.Rows("17:17").RowHeight = 120
